$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("D5").Value = 44242
$ws.Range("J5").Value = 60

# Row 6
$ws.Range("D6").Value = 44242
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 556

# Row 7
$ws.Range("D7").Value = 44756
$ws.Range("J7").Value = 65
$ws.Range("K7").Value = 14000
$ws.Range("L7").Value = 14000
$ws.Range("M7").Value = 14000
$ws.Range("N7").Value = "`$/caja 15 kilos"
$ws.Range("P7").Value = 933
$ws.Range("Q7").Value = 15

# Row 8
$ws.Range("D8").Value = 44756
$ws.Range("J8").Value = 68
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("N8").Value = "`$/caja 15 kilos"
$ws.Range("P8").Value = 800
$ws.Range("Q8").Value = 15

# Row 9
$ws.Range("D9").Value = 44992
$ws.Range("J9").Value = 56
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 13000
$ws.Range("M9").Value = 13000
$ws.Range("N9").Value = "`$/bandeja 18 kilos"
$ws.Range("P9").Value = 722
$ws.Range("Q9").Value = 18

# Row 10
$ws.Range("D10").Value = 44238
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 90
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 13000
$ws.Range("N10").Value = "`$/bandeja 18 kilos"
$ws.Range("P10").Value = 722
$ws.Range("Q10").Value = 18

# Row 11
$ws.Range("D11").Value = 44238
$ws.Range("I11").Value = "Segunda"
$ws.Range("J11").Value = 80
$ws.Range("K11").Value = 11000
$ws.Range("L11").Value = 11000
$ws.Range("M11").Value = 11000
$ws.Range("P11").Value = 611

# Row 12
$ws.Range("D12").Value = 44424
$ws.Range("J12").Value = 75
$ws.Range("K12").Value = 18000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 18000
$ws.Range("P12").Value = 1200

# Row 13
$ws.Range("D13").Value = 44424
$ws.Range("J13").Value = 50

# Row 14
$ws.Range("D14").Value = 44991
$ws.Range("J14").Value = 75

# Row 15
$ws.Range("D15").Value = 44991
$ws.Range("J15").Value = 56
$ws.Range("K15").Value = 9000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 9000
$ws.Range("P15").Value = 500
